$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.440985666666667
$ws.Range("H2").Value = 4.322957000000001
$ws.Range("I2").Value = 0.1098365531732288
$ws.Range("J2").Value = 0.1230162332390494
$ws.Range("M2").Value = 14.349718
$ws.Range("N2").Value = 43.049154
$ws.Range("O2").Value = 0.1016415840981481
$ws.Range("P2").Value = 0.1034081666702025
$ws.Range("Q2").Value = 20.67773795870867
$ws.Range("R2").Value = 186.099641628378
$ws.Range("S2").Value = 0.01116396125640745
$ws.Range("T2").Value = 0.01272088314992413
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.440985666666667
$ws.Range("H3").Value = 4.322957000000001
$ws.Range("I3").Value = 0.1098365531732288
$ws.Range("J3").Value = 0.1230162332390494
$ws.Range("O3").Value = 0.04778708884009916
$ws.Range("P3").Value = 0.04861765281706964
$ws.Range("Q3").Value = 9.721699141278002
$ws.Range("R3").Value = 87.49529227150201
$ws.Range("S3").Value = 0.00524876912437936
$ws.Range("T3").Value = 0.005980760518479764
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.440985666666667
$ws.Range("H4").Value = 4.322957000000001
$ws.Range("I4").Value = 0.1098365531732288
$ws.Range("J4").Value = 0.1230162332390494
$ws.Range("M4").Value = 66.43651233333334
$ws.Range("N4").Value = 199.309537
$ws.Range("O4").Value = 0.4705815372480596
$ws.Range("P4").Value = 0.4787604843769264
$ws.Range("Q4").Value = 95.73406201565659
$ws.Range("R4").Value = 861.6065581409091
$ws.Range("S4").Value = 0.05168705403828625
$ws.Range("T4").Value = 0.05889531141175222
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.440985666666667
$ws.Range("H5").Value = 4.322957000000001
$ws.Range("I5").Value = 0.1098365531732288
$ws.Range("J5").Value = 0.1230162332390494
$ws.Range("M5").Value = 7.2355625
$ws.Range("N5").Value = 14.471125
$ws.Range("O5").Value = 0.05125076564857627
$ws.Range("P5").Value = 0.03476102006337534
$ws.Range("Q5").Value = 10.42634185277084
$ws.Range("R5").Value = 62.55805111662501
$ws.Range("S5").Value = 0.005629207446328535
$ws.Range("T5").Value = 0.004276169751743455
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.440985666666667
$ws.Range("H6").Value = 4.322957000000001
$ws.Range("I6").Value = 0.1098365531732288
$ws.Range("J6").Value = 0.1230162332390494
$ws.Range("M6").Value = 46.41124333333334
$ws.Range("N6").Value = 139.23373
$ws.Range("O6").Value = 0.328739024165117
$ws.Range("P6").Value = 0.3344526760724259
$ws.Range("Q6").Value = 66.87793641551224
$ws.Range("R6").Value = 601.9014277396101
$ws.Range("S6").Value = 0.03610756130782722
$ws.Range("T6").Value = 0.04114310840714977
$ws.Range("I7").Value = 0.5687502547919595
$ws.Range("J7").Value = 0.6369966279614609
$ws.Range("M7").Value = 14.349718
$ws.Range("N7").Value = 43.049154
$ws.Range("O7").Value = 0.1016415840981481
$ws.Range("P7").Value = 0.1034081666702025
$ws.Range("Q7").Value = 107.0724489504773
$ws.Range("R7").Value = 963.6520405542959
$ws.Range("S7").Value = 0.0578086768532801
$ws.Range("T7").Value = 0.06587065347259574
$ws.Range("I8").Value = 0.5687502547919595
$ws.Range("J8").Value = 0.6369966279614609
$ws.Range("O8").Value = 0.04778708884009916
$ws.Range("P8").Value = 0.04861765281706964
$ws.Range("S8").Value = 0.0271789189535724
$ws.Range("T8").Value = 0.03096928090387438
$ws.Range("I9").Value = 0.5687502547919595
$ws.Range("J9").Value = 0.6369966279614609
$ws.Range("M9").Value = 66.43651233333334
$ws.Range("N9").Value = 199.309537
$ws.Range("O9").Value = 0.4705815372480596
$ws.Range("P9").Value = 0.4787604843769264
$ws.Range("Q9").Value = 495.7254264689098
$ws.Range("R9").Value = 4461.528838220188
$ws.Range("S9").Value = 0.2676433692102259
$ws.Range("T9").Value = 0.3049688141492978
$ws.Range("I10").Value = 0.5687502547919595
$ws.Range("J10").Value = 0.6369966279614609
$ws.Range("M10").Value = 7.2355625
$ws.Range("N10").Value = 14.471125
$ws.Range("O10").Value = 0.05125076564857627
$ws.Range("P10").Value = 0.03476102006337534
$ws.Range("Q10").Value = 53.98917221991667
$ws.Range("R10").Value = 323.9350333195
$ws.Range("S10").Value = 0.02914888602091076
$ws.Range("T10").Value = 0.02214265256487078
$ws.Range("I11").Value = 0.5687502547919595
$ws.Range("J11").Value = 0.6369966279614609
$ws.Range("M11").Value = 46.41124333333334
$ws.Range("N11").Value = 139.23373
$ws.Range("O11").Value = 0.328739024165117
$ws.Range("P11").Value = 0.3344526760724259
$ws.Range("Q11").Value = 346.3040515873911
$ws.Range("R11").Value = 3116.73646428652
$ws.Range("S11").Value = 0.1869704037539704
$ws.Range("T11").Value = 0.2130452268708221
$ws.Range("G12").Value = 4.2167365
$ws.Range("H12").Value = 8.433472999999999
$ws.Range("I12").Value = 0.3214131920348118
$ws.Range("J12").Value = 0.2399871387994896
$ws.Range("M12").Value = 14.349718
$ws.Range("N12").Value = 43.049154
$ws.Range("O12").Value = 0.1016415840981481
$ws.Range("P12").Value = 0.1034081666702025
$ws.Range("Q12").Value = 60.508979655307
$ws.Range("R12").Value = 363.053877931842
$ws.Range("S12").Value = 0.03266894598846053
$ws.Range("T12").Value = 0.02481663004768265
$ws.Range("G13").Value = 4.2167365
$ws.Range("H13").Value = 8.433472999999999
$ws.Range("I13").Value = 0.3214131920348118
$ws.Range("J13").Value = 0.2399871387994896
$ws.Range("O13").Value = 0.04778708884009916
$ws.Range("P13").Value = 0.04861765281706964
$ws.Range("Q13").Value = 28.448474234913
$ws.Range("R13").Value = 170.690845409478
$ws.Range("S13").Value = 0.0153594007621474
$ws.Range("T13").Value = 0.01166761139471549
$ws.Range("G14").Value = 4.2167365
$ws.Range("H14").Value = 8.433472999999999
$ws.Range("I14").Value = 0.3214131920348118
$ws.Range("J14").Value = 0.2399871387994896
$ws.Range("M14").Value = 66.43651233333334
$ws.Range("N14").Value = 199.309537
$ws.Range("O14").Value = 0.4705815372480596
$ws.Range("P14").Value = 0.4787604843769264
$ws.Range("Q14").Value = 280.1452664886668
$ws.Range("R14").Value = 1680.871598932001
$ws.Range("S14").Value = 0.1512511139995475
$ws.Range("T14").Value = 0.1148963588158763
$ws.Range("G15").Value = 4.2167365
$ws.Range("H15").Value = 8.433472999999999
$ws.Range("I15").Value = 0.3214131920348118
$ws.Range("J15").Value = 0.2399871387994896
$ws.Range("M15").Value = 7.2355625
$ws.Range("N15").Value = 14.471125
$ws.Range("O15").Value = 0.05125076564857627
$ws.Range("P15").Value = 0.03476102006337534
$ws.Range("Q15").Value = 30.51046049178125
$ws.Range("R15").Value = 122.041841967125
$ws.Range("S15").Value = 0.01647267218133698
$ws.Range("T15").Value = 0.0083421977467611
$ws.Range("G16").Value = 4.2167365
$ws.Range("H16").Value = 8.433472999999999
$ws.Range("I16").Value = 0.3214131920348118
$ws.Range("J16").Value = 0.2399871387994896
$ws.Range("M16").Value = 46.41124333333334
$ws.Range("N16").Value = 139.23373
$ws.Range("O16").Value = 0.328739024165117
$ws.Range("P16").Value = 0.3344526760724259
$ws.Range("Q16").Value = 195.7039837740483
$ws.Range("R16").Value = 1174.22390264429
$ws.Range("S16").Value = 0.1056610591033194
$ws.Range("T16").Value = 0.08026434079445401
